# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across several sheets with latest market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 950.2
$ws.Range("I2").Value = 950.2
$ws.Range("K2").Value = 950.2
$ws.Range("M2").Value = -837.2

$ws.Range("H5").Value = 76
$ws.Range("I5").Value = 76
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 76
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 39
$ws.Range("N5").ClearContents()

$ws.Range("H17").Value = 2000
$ws.Range("J17").Value = 2000
$ws.Range("L17").Value = 6000
$ws.Range("N17").Value = -6336

$ws.Range("H55").Value = 274.75
$ws.Range("I55").Value = 349.66666
$ws.Range("J55").Value = 50
$ws.Range("K55").Value = 349.66666
$ws.Range("L55").Value = 50
$ws.Range("M55").Value = -135.66666
$ws.Range("N55").Value = -478

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6003.1665
$ws.Range("I45").Value = 6203.8
$ws.Range("K45").Value = 6203.8
$ws.Range("M45").Value = -5826.8

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H88").Value = 893
$ws.Range("I88").Value = 998.5
$ws.Range("J88").Value = 787.5
$ws.Range("K88").Value = 998.5
$ws.Range("L88").Value = 787.5
$ws.Range("M88").Value = -592.5
$ws.Range("N88").Value = -1599.5

$ws.Range("H91").Value = 893
$ws.Range("I91").Value = 998.5
$ws.Range("J91").Value = 787.5
$ws.Range("K91").Value = 998.5
$ws.Range("L91").Value = 787.5
$ws.Range("M91").Value = 405.5
$ws.Range("N91").Value = -3595.5

$ws.Range("H101").Value = 35000
$ws.Range("J101").Value = 35000
$ws.Range("L101").Value = 35000
$ws.Range("N101").Value = -41490

$ws.Range("H122").Value = 3192.6667
$ws.Range("I122").Value = 3192.6667
$ws.Range("K122").Value = 9578.000100000001
$ws.Range("M122").Value = -7128.000100000001

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 3470.4546
$ws.Range("I132").Value = 3017.5
$ws.Range("K132").Value = 9052.5
$ws.Range("M132").Value = -6522.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 191
$ws.Range("I22").Value = 194.4
$ws.Range("J22").Value = 174
$ws.Range("K22").Value = 194.4
$ws.Range("L22").Value = 174
$ws.Range("M22").Value = -21.40000000000001
$ws.Range("N22").Value = -520

$ws.Range("H94").Value = 1666.4445
$ws.Range("J94").Value = 150
$ws.Range("L94").Value = 150
$ws.Range("N94").Value = -1052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 282.41666
$ws.Range("I7").Value = 332.2
$ws.Range("J7").Value = 33.5
$ws.Range("K7").Value = 332.2
$ws.Range("L7").Value = 33.5
$ws.Range("M7").Value = -219.2
$ws.Range("N7").Value = -259.5

$ws.Range("H22").Value = 485.85715
$ws.Range("I22").Value = 280.2
$ws.Range("K22").Value = 280.2
$ws.Range("M22").Value = 69.80000000000001

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H58").Value = 9584.666999999999
$ws.Range("I58").Value = 7377.5
$ws.Range("J58").Value = 13999
$ws.Range("K58").Value = 7377.5
$ws.Range("L58").Value = 13999
$ws.Range("M58").Value = -7174.5
$ws.Range("N58").Value = -14405

$ws.Range("H132").Value = 8858.6
$ws.Range("I132").Value = 4764.3335
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 14293.0005
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -11763.0005
$ws.Range("N132").Value = -50060

$ws.Range("H136").Value = 9584.666999999999
$ws.Range("I136").Value = 7377.5
$ws.Range("J136").Value = 13999
$ws.Range("K136").Value = 22132.5
$ws.Range("L136").Value = 41997
$ws.Range("M136").Value = -19582.5
$ws.Range("N136").Value = -47097

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 706.5625
$ws.Range("I2").Value = 375.22223
$ws.Range("J2").Value = 2495.8
$ws.Range("K2").Value = 2251.33338
$ws.Range("L2").Value = 14974.8
$ws.Range("M2").Value = -2138.33338
$ws.Range("N2").Value = -15200.8

$ws.Range("H12").Value = 42.076923
$ws.Range("J12").Value = 35.666668
$ws.Range("L12").Value = 107.000004
$ws.Range("N12").Value = -453.000004

$ws.Range("H50").Value = 199
$ws.Range("I50").Value = 199
$ws.Range("K50").Value = 597
$ws.Range("M50").Value = -116

$ws.Range("H53").Value = 199
$ws.Range("I53").Value = 199
$ws.Range("K53").Value = 597
$ws.Range("M53").Value = -116

$ws.Range("H80").Value = 1497.5
$ws.Range("J80").Value = 1497.5
$ws.Range("L80").Value = 4492.5
$ws.Range("N80").Value = -6364.5

$ws.Range("H83").Value = 1497.5
$ws.Range("J83").Value = 1497.5
$ws.Range("L83").Value = 13477.5
$ws.Range("N83").Value = -22837.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3499.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3499.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3499.5
$ws.Range("N80").Value = -5495.5
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 3499.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3499.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 17497.5
$ws.Range("N83").Value = -27481.5
$ws.Range("M83").ClearContents()

$ws.Range("H107").Value = 322.5
$ws.Range("I107").Value = 195
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 195
$ws.Range("L107").Value = 450
$ws.Range("M107").Value = 1725
$ws.Range("N107").Value = -4290

$ws.Range("H118").Value = 39999
$ws.Range("J118").Value = 39999
$ws.Range("L118").Value = 39999
$ws.Range("N118").Value = -43313

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 933.3333
$ws.Range("I22").Value = 933.3333
$ws.Range("K22").Value = 933.3333
$ws.Range("M22").Value = -638.3333

$ws.Range("H27").Value = 933.3333
$ws.Range("I27").Value = 933.3333
$ws.Range("K27").Value = 933.3333
$ws.Range("M27").Value = -826.3333

$ws.Range("H68").Value = 3200
$ws.Range("I68").Value = 3200
$ws.Range("K68").Value = 3200
$ws.Range("M68").Value = -2451

$ws.Range("H71").Value = 3200
$ws.Range("I71").Value = 3200
$ws.Range("K71").Value = 16000
$ws.Range("M71").Value = -12256

$ws.Range("H122").Value = 3700
$ws.Range("I122").Value = 3700
$ws.Range("K122").Value = 11100
$ws.Range("M122").Value = -8650

$ws.Range("H132").Value = 6705.4287
$ws.Range("I132").Value = 6705.4287
$ws.Range("K132").Value = 20116.2861
$ws.Range("M132").Value = -17586.2861

